$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4278.3335
$ws.Range("I113").Value = 2990
$ws.Range("J113").Value = 5098.1816
$ws.Range("K113").Value = 2990
$ws.Range("L113").Value = 5098.1816
$ws.Range("M113").Value = 264
$ws.Range("N113").Value = -11606.1816

$ws.Range("H125").Value = 45454772
$ws.Range("I125").Value = 350
$ws.Range("J125").Value = 71428730
$ws.Range("K125").Value = 3150
$ws.Range("L125").Value = 642858570
$ws.Range("M125").Value = -690
$ws.Range("N125").Value = -642863490

$ws.Range("H132").Value = 3042.4133
$ws.Range("I132").Value = 1470.74
$ws.Range("J132").Value = 6185.76
$ws.Range("K132").Value = 4412.22
$ws.Range("L132").Value = 18557.28
$ws.Range("M132").Value = -1882.22
$ws.Range("N132").Value = -23617.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3920.19
$ws.Range("I32").Value = 2537.4658
$ws.Range("J32").Value = 14060.167
$ws.Range("K32").Value = 2537.4658
$ws.Range("L32").Value = 14060.167
$ws.Range("M32").Value = -2250.4658
$ws.Range("N32").Value = -14634.167

$ws.Range("H61").Value = 264821.56
$ws.Range("I61").Value = 246074.33
$ws.Range("J61").Value = 288925.16
$ws.Range("K61").Value = 246074.33
$ws.Range("L61").Value = 288925.16
$ws.Range("M61").Value = -245862.33
$ws.Range("N61").Value = -289349.16

$ws.Range("H101").Value = 31995
$ws.Range("J101").Value = 31995
$ws.Range("L101").Value = 31995
$ws.Range("N101").Value = -38485

$ws.Range("H102").Value = 2008.2444
$ws.Range("I102").Value = 1417.5
$ws.Range("J102").Value = 3462.3845
$ws.Range("K102").Value = 1417.5
$ws.Range("L102").Value = 3462.3845
$ws.Range("M102").Value = 204.5
$ws.Range("N102").Value = -6706.3845

$ws.Range("H132").Value = 2087.2173
$ws.Range("I132").Value = 1941.5933
$ws.Range("K132").Value = 5824.7799
$ws.Range("M132").Value = -3294.7799

$ws.Range("H135").Value = 20238.9
$ws.Range("J135").Value = 20238.9
$ws.Range("L135").Value = 20238.9
$ws.Range("N135").Value = -30378.9

$ws.Range("H136").Value = 264821.56
$ws.Range("I136").Value = 246074.33
$ws.Range("J136").Value = 288925.16
$ws.Range("K136").Value = 738222.99
$ws.Range("L136").Value = 866775.48
$ws.Range("M136").Value = -735672.99
$ws.Range("N136").Value = -871875.48

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1766.7037
$ws.Range("I105").Value = 1556.6666
$ws.Range("J105").Value = 2501.8333
$ws.Range("K105").Value = 1556.6666
$ws.Range("L105").Value = 2501.8333
$ws.Range("M105").Value = 190.3334
$ws.Range("N105").Value = -5995.8333

$ws.Range("H135").Value = 22825.715
$ws.Range("J135").Value = 21956
$ws.Range("L135").Value = 21956
$ws.Range("N135").Value = -32096

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3249.157
$ws.Range("I31").Value = 2296.8064
$ws.Range("J31").Value = 4725.3
$ws.Range("K31").Value = 2296.8064
$ws.Range("L31").Value = 4725.3
$ws.Range("M31").Value = -2001.8064
$ws.Range("N31").Value = -5315.3

$ws.Range("H34").Value = 3249.157
$ws.Range("I34").Value = 2296.8064
$ws.Range("J34").Value = 4725.3
$ws.Range("K34").Value = 2296.8064
$ws.Range("L34").Value = 4725.3
$ws.Range("M34").Value = -2094.8064
$ws.Range("N34").Value = -5129.3

$ws.Range("H99").Value = 37271.605
$ws.Range("I99").Value = 72699.36
$ws.Range("J99").Value = 1843.8572
$ws.Range("K99").Value = 72699.36
$ws.Range("L99").Value = 1843.8572
$ws.Range("M99").Value = -71201.36
$ws.Range("N99").Value = -4839.8572

$ws.Range("H105").Value = 668.75
$ws.Range("J105").Value = 875
$ws.Range("L105").Value = 875
$ws.Range("N105").Value = -4369

$ws.Range("H126").Value = 37271.605
$ws.Range("I126").Value = 72699.36
$ws.Range("J126").Value = 1843.8572
$ws.Range("K126").Value = 218098.08
$ws.Range("L126").Value = 5531.571599999999
$ws.Range("M126").Value = -215628.08
$ws.Range("N126").Value = -10471.5716

$ws.Range("H132").Value = 1644.8392
$ws.Range("I132").Value = 1054.3513
$ws.Range("J132").Value = 2794.7368
$ws.Range("K132").Value = 3163.0539
$ws.Range("L132").Value = 8384.2104
$ws.Range("M132").Value = -633.0538999999999
$ws.Range("N132").Value = -13444.2104

$ws.Range("H134").Value = 1494.7906
$ws.Range("I134").Value = 1010.4717
$ws.Range("J134").Value = 2272.6365
$ws.Range("K134").Value = 3031.4151
$ws.Range("L134").Value = 6817.9095
$ws.Range("M134").Value = -496.4151000000002
$ws.Range("N134").Value = -11887.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 7550
$ws.Range("J102").Value = 7611.1113
$ws.Range("L102").Value = 22833.3339
$ws.Range("N102").Value = -27701.3339

$ws.Range("H131").Value = 1124.2766
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1124.2766
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 3372.8298
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -13452.8298

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4618.9473
$ws.Range("I102").Value = 2983.0667
$ws.Range("J102").Value = 10753.5
$ws.Range("K102").Value = 2983.0667
$ws.Range("L102").Value = 10753.5
$ws.Range("M102").Value = -1361.0667
$ws.Range("N102").Value = -13997.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1633.1708
$ws.Range("I7").Value = 1290.909
$ws.Range("J7").Value = 2029.4736
$ws.Range("K7").Value = 1290.909
$ws.Range("L7").Value = 2029.4736
$ws.Range("M7").Value = -1178.909
$ws.Range("N7").Value = -2253.4736

$ws.Range("H61").Value = 2363
$ws.Range("I61").Value = 2272
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2272
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2070
$ws.Range("N61").Value = -3404

$ws.Range("H100").Value = 83339280
$ws.Range("I100").Value = 11820
$ws.Range("J100").Value = 142858900
$ws.Range("K100").Value = 11820
$ws.Range("L100").Value = 142858900
$ws.Range("M100").Value = -11279
$ws.Range("N100").Value = -142859982

$ws.Range("H113").Value = 2363
$ws.Range("I113").Value = 2272
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2272
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -102
$ws.Range("N113").Value = -7340

$ws.Range("H126").Value = 1633.1708
$ws.Range("I126").Value = 1290.909
$ws.Range("J126").Value = 2029.4736
$ws.Range("K126").Value = 3872.727
$ws.Range("L126").Value = 6088.4208
$ws.Range("M126").Value = -1402.727
$ws.Range("N126").Value = -11028.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4297.3335
$ws.Range("I4").Value = 2350
$ws.Range("J4").Value = 4853.7144
$ws.Range("K4").Value = 2350
$ws.Range("L4").Value = 4853.7144
$ws.Range("M4").Value = -2237
$ws.Range("N4").Value = -5079.7144

$ws.Range("H132").Value = 1811.3707
$ws.Range("I132").Value = 1304.8292
$ws.Range("K132").Value = 3914.487599999999
$ws.Range("M132").Value = -1384.487599999999
